{"js": "// Remove the \"\u2705\" (checkmark emoji) runs from the \"M\u1ee9c \u0111\u1ed9 ho\u00e0n th\u00e0nh\"\n// column of the progress table, and drop the now-unneeded leading\n// space in front of \"100%\" so the cell just reads \"100%\".\n\nconst body = context.document.body;\n\n// --- Step 1: delete every \"\u2705\" run in the document -----------------\nconst emojiResults = body.search(\"\u2705\", { matchCase: true });\nemojiResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < emojiResults.items.length; i++) {\n  emojiResults.items[i].delete();\n}\nawait context.sync();\n\n// --- Step 2: strip the leading space that used to separate the\n//     emoji from \"100%\" (\" 100%\" -> \"100%\") without disturbing the\n//     \"100%\" run itself. --------------------------------------------\nconst spaceResults = body.search(\" 100%\", { matchCase: true });\nspaceResults.load(\"items\");\nawait context.sync();\n\n// Split each match into its individual \" \" / \"100%\" pieces so we can\n// delete just the leading space character and leave \"100%\" intact.\nconst splitRanges = [];\nfor (let i = 0; i < spaceResults.items.length; i++) {\n  const pieces = spaceResults.items[i].getTextRanges([\" \"], false);\n  pieces.load(\"items\");\n  splitRanges.push(pieces);\n}\nawait context.sync();\n\nfor (let i = 0; i < splitRanges.length; i++) {\n  const pieces = splitRanges[i].items;\n  if (pieces.length > 0) {\n    pieces[0].delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Step 1: remove every \"\u2705\" (checkmark emoji) run from the\n# \"M\u1ee9c \u0111\u1ed9 ho\u00e0n th\u00e0nh\" column of the progress table.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"\u2705\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# Step 2: the emoji used to be followed by a space before \"100%\"; now\n# that the emoji is gone, drop that now-orphaned leading space so the\n# cell just reads \"100%\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\" 100%\", $false, $false, $false, $false, $false, $true, 1, $false, \"100%\", 2) | Out-Null\n"}
